# Auto-generated Excel COM-interop edit script
# Applies numeric updates to the Bahamut_Profits leve-crafting tables
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) per the scheduled-runner price refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 74
$ws.Range("H74").Value = 4723420.5
$ws.Range("I74").Value = 8654770
$ws.Range("J74").Value = 5800
$ws.Range("K74").Value = 8654770
$ws.Range("L74").Value = 5800
$ws.Range("M74").Value = -8653834
$ws.Range("N74").Value = -7672
# row 77
$ws.Range("H77").Value = 4723420.5
$ws.Range("I77").Value = 8654770
$ws.Range("J77").Value = 5800
$ws.Range("K77").Value = 43273850
$ws.Range("L77").Value = 29000
$ws.Range("M77").Value = -43269170
$ws.Range("N77").Value = -38360
# row 107
$ws.Range("H107").Value = 769328.25
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
# row 113
$ws.Range("H113").Value = 3340.9285
$ws.Range("I113").Value = 3388.125
$ws.Range("J113").Value = 3278
$ws.Range("K113").Value = 3388.125
$ws.Range("L113").Value = 3278
$ws.Range("M113").Value = -134.125
$ws.Range("N113").Value = -9786
# row 129
$ws.Range("H129").Value = 1099.4857
$ws.Range("I129").Value = 493
$ws.Range("J129").Value = 1177.742
$ws.Range("K129").Value = 1479
$ws.Range("L129").Value = 3533.226
$ws.Range("M129").Value = 3521
$ws.Range("N129").Value = -13533.226
# row 132
$ws.Range("H132").Value = 264945.84
$ws.Range("I132").Value = 1840.8387
$ws.Range("J132").Value = 1430125.1
$ws.Range("K132").Value = 5522.5161
$ws.Range("L132").Value = 4290375.300000001
$ws.Range("M132").Value = -2992.5161
$ws.Range("N132").Value = -4295435.300000001

$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 3181.9167
$ws.Range("I2").Value = 2874.5715
$ws.Range("J2").Value = 5333.3335
$ws.Range("K2").Value = 2874.5715
$ws.Range("L2").Value = 5333.3335
$ws.Range("M2").Value = -2761.5715
$ws.Range("N2").Value = -5559.3335
# row 5
$ws.Range("H5").Value = 109.5
$ws.Range("I5").Value = 114.5
$ws.Range("J5").Value = 99.5
$ws.Range("K5").Value = 114.5
$ws.Range("L5").Value = 99.5
$ws.Range("M5").Value = -2.5
$ws.Range("N5").Value = -323.5
# row 32
$ws.Range("H32").Value = 4361.547
$ws.Range("I32").Value = 3588.6128
$ws.Range("J32").Value = 8047.846
$ws.Range("K32").Value = 3588.6128
$ws.Range("L32").Value = 8047.846
$ws.Range("M32").Value = -3301.6128
$ws.Range("N32").Value = -8621.846
# row 45
$ws.Range("H45").Value = 1926.5555
$ws.Range("I45").Value = 1999.5
$ws.Range("J45").Value = 1905.7142
$ws.Range("K45").Value = 1999.5
$ws.Range("L45").Value = 1905.7142
$ws.Range("M45").Value = -1622.5
$ws.Range("N45").Value = -2659.7142
# row 116
$ws.Range("H116").Value = 3181.9167
$ws.Range("I116").Value = 2874.5715
$ws.Range("J116").Value = 5333.3335
$ws.Range("K116").Value = 2874.5715
$ws.Range("L116").Value = 5333.3335
$ws.Range("M116").Value = -580.5715
$ws.Range("N116").Value = -9921.333500000001
# row 122
$ws.Range("H122").Value = 1271.8334
$ws.Range("I122").Value = 1271.8334
$ws.Range("K122").Value = 3815.5002
$ws.Range("M122").Value = -1365.5002
# row 132
$ws.Range("H132").Value = 2118.606
$ws.Range("I132").Value = 1565.48
$ws.Range("J132").Value = 3847.125
$ws.Range("K132").Value = 4696.440000000001
$ws.Range("L132").Value = 11541.375
$ws.Range("M132").Value = -2166.440000000001
$ws.Range("N132").Value = -16601.375

$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 3181.9167
$ws.Range("I3").Value = 2874.5715
$ws.Range("J3").Value = 5333.3335
$ws.Range("K3").Value = 2874.5715
$ws.Range("L3").Value = 5333.3335
$ws.Range("M3").Value = -2760.5715
$ws.Range("N3").Value = -5561.3335
# row 4
$ws.Range("H4").Value = 109.5
$ws.Range("I4").Value = 114.5
$ws.Range("J4").Value = 99.5
$ws.Range("K4").Value = 114.5
$ws.Range("L4").Value = 99.5
$ws.Range("M4").Value = 0.5
$ws.Range("N4").Value = -329.5
# row 134
$ws.Range("H134").Value = 28548.648
$ws.Range("I134").Value = 1432.4517
$ws.Range("J134").Value = 168649
$ws.Range("K134").Value = 4297.355100000001
$ws.Range("L134").Value = 505947
$ws.Range("M134").Value = -1762.355100000001
$ws.Range("N134").Value = -511017

$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 1303.0454
$ws.Range("I16").Value = 1271.2632
$ws.Range("J16").Value = 1504.3334
$ws.Range("K16").Value = 1271.2632
$ws.Range("L16").Value = 1504.3334
$ws.Range("M16").Value = -984.2632000000001
$ws.Range("N16").Value = -2078.3334
# row 31
$ws.Range("H31").Value = 2195.7
$ws.Range("I31").Value = 2191.8215
$ws.Range("J31").Value = 2250
$ws.Range("K31").Value = 2191.8215
$ws.Range("L31").Value = 2250
$ws.Range("M31").Value = -1896.8215
$ws.Range("N31").Value = -2840
# row 34
$ws.Range("H34").Value = 2195.7
$ws.Range("I34").Value = 2191.8215
$ws.Range("J34").Value = 2250
$ws.Range("K34").Value = 2191.8215
$ws.Range("L34").Value = 2250
$ws.Range("M34").Value = -1989.8215
$ws.Range("N34").Value = -2654
# row 86
$ws.Range("H86").Value = 3666.6667
$ws.Range("K86").Value = 3750
$ws.Range("M86").Value = -2627
# row 89
$ws.Range("H89").Value = 3666.6667
$ws.Range("K89").Value = 18750
$ws.Range("M89").Value = -13134
# row 113
$ws.Range("H113").Value = 1303.0454
$ws.Range("I113").Value = 1271.2632
$ws.Range("J113").Value = 1504.3334
$ws.Range("K113").Value = 1271.2632
$ws.Range("L113").Value = 1504.3334
$ws.Range("M113").Value = 898.7367999999999
$ws.Range("N113").Value = -5844.3334
# row 132
$ws.Range("H132").Value = 2258.0286
$ws.Range("I132").Value = 1522.2916
$ws.Range("K132").Value = 4566.8748
$ws.Range("M132").Value = -2036.8748
# row 134
$ws.Range("H134").Value = 1831.721
$ws.Range("I134").Value = 1306.7693
$ws.Range("J134").Value = 6950
$ws.Range("K134").Value = 3920.3079
$ws.Range("L134").Value = 20850
$ws.Range("M134").Value = -1385.3079
$ws.Range("N134").Value = -25920

$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 1423.4736
$ws.Range("J5").Value = 3745
$ws.Range("L5").Value = 11235
$ws.Range("N5").Value = -11459
# row 127
$ws.Range("H127").Value = 1402.875
$ws.Range("J127").Value = 1402.875
$ws.Range("L127").Value = 4208.625
$ws.Range("N127").Value = -14128.625
# row 135
$ws.Range("H135").Value = 1423.4736
$ws.Range("J135").Value = 3745
$ws.Range("L135").Value = 33705
$ws.Range("N135").Value = -38775
# row 138
$ws.Range("H138").Value = 3060.8125
$ws.Range("I138").Value = 2444
$ws.Range("J138").Value = 4088.8333
$ws.Range("K138").Value = 7332
$ws.Range("L138").Value = 12266.4999
$ws.Range("M138").Value = -2192
$ws.Range("N138").Value = -22546.4999

$ws = $wb.Worksheets.Item("GSM")
# row 70
$ws.Range("H70").Value = 4541.4707
$ws.Range("I70").Value = 4080.3333
$ws.Range("J70").Value = 8000
$ws.Range("K70").Value = 4080.3333
$ws.Range("L70").Value = 8000
$ws.Range("M70").Value = -3810.3333
$ws.Range("N70").Value = -8540
# row 73
$ws.Range("H73").Value = 4541.4707
$ws.Range("I73").Value = 4080.3333
$ws.Range("J73").Value = 8000
$ws.Range("K73").Value = 4080.3333
$ws.Range("L73").Value = 8000
$ws.Range("M73").Value = -3144.3333
$ws.Range("N73").Value = -9872
# row 100
$ws.Range("H100").Value = 19190
$ws.Range("J100").Value = 19190
$ws.Range("L100").Value = 19190
$ws.Range("N100").Value = -21354
# row 132
$ws.Range("H132").Value = 1767.3103
$ws.Range("I132").Value = 1488.1177
$ws.Range("J132").Value = 2162.8333
$ws.Range("K132").Value = 4464.3531
$ws.Range("L132").Value = 6488.499899999999
$ws.Range("M132").Value = -1934.3531
$ws.Range("N132").Value = -11548.4999

$ws = $wb.Worksheets.Item("LTW")
# row 6
$ws.Range("H6").Value = 20000
$ws.Range("J6").Value = 20000
$ws.Range("L6").Value = 20000
$ws.Range("N6").Value = -20224
# row 136
$ws.Range("H136").Value = 2325
$ws.Range("I136").Value = 1301.4062
$ws.Range("J136").Value = 7784.1665
$ws.Range("K136").Value = 3904.2186
$ws.Range("L136").Value = 23352.4995
$ws.Range("M136").Value = -1354.2186
$ws.Range("N136").Value = -28452.4995

$ws = $wb.Worksheets.Item("WVR")
# row 107
$ws.Range("H107").Value = 610.7778
$ws.Range("I107").Value = 699.0769
$ws.Range("J107").Value = 381.2
$ws.Range("K107").Value = 2097.2307
$ws.Range("L107").Value = 1143.6
$ws.Range("M107").Value = -177.2307000000001
$ws.Range("N107").Value = -4983.6
# row 113
$ws.Range("H113").Value = 207.0625
$ws.Range("I113").Value = 196.03703
$ws.Range("K113").Value = 588.11109
$ws.Range("M113").Value = 1581.88891
# row 123
$ws.Range("H123").Value = 30265.385
$ws.Range("J123").Value = 30265.385
$ws.Range("L123").Value = 30265.385
$ws.Range("N123").Value = -40065.38499999999
# row 132
$ws.Range("H132").Value = 1519.4286
$ws.Range("I132").Value = 1202.1765
$ws.Range("J132").Value = 2009.7273
$ws.Range("K132").Value = 3606.5295
$ws.Range("L132").Value = 6029.1819
$ws.Range("M132").Value = -1076.5295
$ws.Range("N132").Value = -11089.1819
# row 135
$ws.Range("H135").Value = 60255
$ws.Range("J135").Value = 60255
$ws.Range("L135").Value = 60255
$ws.Range("N135").Value = -70395
# row 136
$ws.Range("H136").Value = 1756.4791
$ws.Range("I136").Value = 1696.697
$ws.Range("J136").Value = 1888
$ws.Range("K136").Value = 5090.090999999999
$ws.Range("L136").Value = 5664
$ws.Range("M136").Value = -2540.090999999999
$ws.Range("N136").Value = -10764

# ALC row 107: N107 column (LeveProfitHQ) no longer applies - clear it
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N107").ClearContents()
